$wb = $excel.ActiveWorkbook

# --- Sheet: Change Management Overview ---
$ws1 = $wb.Worksheets.Item("Change Management Overview")

$ws1.Range("B6").Value = "Enterprise Cloud Infrastructure Migration"

# materialize an empty row 13 between the existing row 12 and row 14
$ws1.Rows.Item(13).OutlineLevel = 0

$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new IT systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in IT technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for IT transformation"

# materialize an empty row 21 between the existing row 20 and row 22
$ws1.Rows.Item(21).OutlineLevel = 0

# --- Sheet: Change Impact Assessment ---
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

# materialize an empty row 2 between the existing row 1 and row 3
$ws2.Rows.Item(2).OutlineLevel = 0

$ws2.Range("A4").Value = "IT Managers"
$ws2.Range("G4").Value = "IT automation"
$ws2.Range("A5").Value = "System Administrators"

# --- Sheet: Change Activities ---
$ws3 = $wb.Worksheets.Item("Change Activities")

# materialize an empty row 2 between the existing row 1 and row 3
$ws3.Rows.Item(2).OutlineLevel = 0
